$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 1. Rename the loan product (productname) and set the numeric shortname.
#    Clear first so the freed shared-string slots get reclaimed/reused in the
#    same order the target workbook uses them.
# ---------------------------------------------------------------------------
$ws1.Range("B1").ClearContents()
$ws1.Range("B3").ClearContents()
$ws2.Range("B1").ClearContents()

$ws1.Range("B1").Value = "438-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"
$ws2.Range("B1").Value = "438-RBI-EI-DB-SAR-REC-NON-RNI-CTPD-DL-MD-TR-1-EarlyRePayment"
$ws1.Range("B3").Value = 438

# ---------------------------------------------------------------------------
# 2. nominalinterestratedefault changes from 12 to 1.
# ---------------------------------------------------------------------------
$ws1.Range("B11").Value = 1

# ---------------------------------------------------------------------------
# 3. Insert the new "preclosureinterestcalculationrule" row right after the
#    "recalculateinterest" row (row 21), pushing the rest down by one.
# ---------------------------------------------------------------------------
$ws1.Rows.Item(22).Insert()
$ws1.Range("A22").Value = "preclosureinterestcalculationrule"
$ws1.Range("B22").Value = "Calculate till pre closure date"

# ---------------------------------------------------------------------------
# 4. Append the new fund-source / GL-account mapping rows (31-42). Copy the
#    formatting (s="3" on column A, s="1" on column B) from an existing row
#    that already carries that style pairing before writing the new values.
#    The values (column B) were typed in before the labels (column A), so
#    write them in that order to match the resulting shared-string table.
# ---------------------------------------------------------------------------
$ws1.Range("A31:B42").ClearContents()
$ws1.Range("A13:B13").Copy()
$ws1.Range("A31:B42").PasteSpecial(-4122)

$ws1.Range("B31").Value = "Cash"
$ws1.Range("B32").Value = "Loan portfolio "
$ws1.Range("B33").Value = "Interest Receivable "
$ws1.Range("B34").Value = "Penalties Receivable "
$ws1.Range("B35").Value = "Transfer in Suspence "
$ws1.Range("B36").Value = "Fees Receivable"
$ws1.Range("B37").Value = "Income from interest"
$ws1.Range("B38").Value = "Income from penalties"
$ws1.Range("B39").Value = "Income from fees"
$ws1.Range("B40").Value = "Income from recovery repayments"
$ws1.Range("B41").Value = "Losses Writtenoff "
$ws1.Range("B42").Value = "Overpayment Liability"

$ws1.Range("A31").Value = "fundsource"
$ws1.Range("A32").Value = "loanprotfolio"
$ws1.Range("A33").Value = "interestreceivable"
$ws1.Range("A34").Value = "penaltiesreceivable"
$ws1.Range("A35").Value = "transferinsuspense"
$ws1.Range("A36").Value = "feesreceivable"
$ws1.Range("A37").Value = "incomefrominterest"
$ws1.Range("A38").Value = "incomefrompenalties"
$ws1.Range("A39").Value = "incomefromfees"
$ws1.Range("A40").Value = "incomefromrecoveryrepayments"
$ws1.Range("A41").Value = "loseswrittenoff"
$ws1.Range("A42").Value = "overpaymentliability"

# ---------------------------------------------------------------------------
# 5. Restore selections (sheet1 scrolled/selected at A44, sheet2 at B29) while
#    keeping sheet2 the active/tab-selected sheet, matching the workbook's
#    activeTab.
# ---------------------------------------------------------------------------
$ws1.Select()
$ws1.Range("A44").Select()
$ws2.Select()
$ws2.Range("B29").Select()
